$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D-column (Price) values as text to preserve exact formatting (avoid Excel numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.205.22"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.17"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.54"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.82"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2849"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06469"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.77"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07690"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.60"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.47"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.079"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.6827"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "268.92"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.193.54"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.34"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007527"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.101.05"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.191"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.106"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.322"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.894"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.374"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09803"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.240"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.979"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04691"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6843"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.741"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.364"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8377"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.886"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.208"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.943"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "918.97"

# Set B, C, E column values (plain text / percent strings, safe to assign directly)
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  -5.68%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -0.50%  "
